$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp title
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 12:20"

# Row 5: España -> España
$ws.Range("F5").Value = 6416

# Row 7: Alemania -> Alemania
$ws.Range("B7").Value = 85063
$ws.Range("C7").Value = 269
$ws.Range("E7").Value = 61512
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 1111

# Row 10: Iran -> Iran
$ws.Range("B10").Value = 53183
$ws.Range("C10").Value = 2715
$ws.Range("D10").Value = 17935
$ws.Range("E10").Value = 31954
$ws.Range("F10").Value = 4035
$ws.Range("G10").Value = 134
$ws.Range("H10").Value = 3294

# Row 17: Austria -> Austria
$ws.Range("B17").Value = 11251
$ws.Range("C17").Value = 122
$ws.Range("E17").Value = 9061
$ws.Range("F17").Value = 245
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 168

# Row 24: Noruega -> Noruega
$ws.Range("E24").Value = 5169
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 54

# Row 31: Ecuador -> Rumania
$ws.Range("A31").Value = "Rumania"
$ws.Range("B31").Value = 3183
$ws.Range("C31").Value = 445
$ws.Range("D31").Value = 267
$ws.Range("E31").Value = 2800
$ws.Range("F31").Value = 78
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 116

# Row 32: Polonia -> Ecuador
$ws.Range("A32").Value = "Ecuador"
$ws.Range("B32").Value = 3163
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 65
$ws.Range("E32").Value = 2978
$ws.Range("F32").Value = 100
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 120

# Row 33: Filipinas -> Polonia
$ws.Range("A33").Value = "Polonia"
$ws.Range("B33").Value = 3149
$ws.Range("C33").Value = 203
$ws.Range("D33").Value = 56
$ws.Range("E33").Value = 3034
$ws.Range("F33").Value = 50
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 59

# Row 34: Rumania -> Filipinas
$ws.Range("A34").Value = "Filipinas"
$ws.Range("B34").Value = 3018
$ws.Range("C34").Value = 385
$ws.Range("D34").Value = 52
$ws.Range("E34").Value = 2830
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 29
$ws.Range("H34").Value = 136

# Row 72: Libano -> Libano
$ws.Range("F72").Value = 2

# Row 74: Letonia -> Letonia
$ws.Range("E74").Value = 461
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 1

# Row 78: Azerbaiyan -> Eslovaquia
$ws.Range("A78").Value = "Eslovaquia"
$ws.Range("B78").Value = 450
$ws.Range("C78").Value = 24
$ws.Range("D78").Value = 5
$ws.Range("E78").Value = 444
$ws.Range("F78").Value = 3
$ws.Range("H78").Value = 1

# Row 79: Principado de Andorra -> Azerbaiyan
$ws.Range("A79").Value = "Azerbaiyan"
$ws.Range("B79").Value = 443
$ws.Range("C79").Value = 43
$ws.Range("D79").Value = 32
$ws.Range("E79").Value = 406
$ws.Range("F79").Value = 7
$ws.Range("H79").Value = 5

# Row 80: Eslovaquia -> Principado de Andorra
$ws.Range("A80").Value = "Principado de Andorra"
$ws.Range("B80").Value = 428
$ws.Range("D80").Value = 10
$ws.Range("E80").Value = 403
$ws.Range("F80").Value = 12
$ws.Range("H80").Value = 15

# Row 90: Jordania -> Albania
$ws.Range("A90").Value = "Albania"
$ws.Range("B90").Value = 304
$ws.Range("C90").Value = 27
$ws.Range("D90").Value = 89
$ws.Range("E90").Value = 199
$ws.Range("F90").Value = 7
$ws.Range("H90").Value = 16

# Row 91: Burkina Faso -> Jordania
$ws.Range("A91").Value = "Jordania"
$ws.Range("B91").Value = 299
$ws.Range("D91").Value = 45
$ws.Range("E91").Value = 249
$ws.Range("F91").Value = 5
$ws.Range("H91").Value = 5

# Row 92: Albania -> Burkina Faso
$ws.Range("A92").Value = "Burkina Faso"
$ws.Range("B92").Value = 288
$ws.Range("D92").Value = 50
$ws.Range("E92").Value = 222
$ws.Range("F92").Value = 0
